# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values that look like plain decimal numbers (single "." and only digits)
# must be forced to Text format first, otherwise Excel auto-converts the typed-in
# string into a numeric value (exactly like typing it into a cell would).
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D10", "D11", "D12",
    "D14", "D15", "D16", "D18", "D19", "D22", "D23", "D24",
    "D25", "D26", "D28", "D29", "D30", "D31", "D33", "D34",
    "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42",
    "D45", "D46", "D47", "D48", "D49", "D50"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '30.080.91'
$ws.Range("E2").Value = '  -3.36%  '
$ws.Range("D3").Value = '1.857.84'
$ws.Range("E3").Value = '  -4.18%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").Value = '232.94'
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").Value = '0.4660'
$ws.Range("E7").Value = '  -2.98%  '
$ws.Range("D8").Value = '0.2804'
$ws.Range("E8").Value = '  -3.26%  '
$ws.Range("E9").Value = '  -4.19%  '
$ws.Range("D10").Value = '19.45'
$ws.Range("E10").Value = '  -1.83%  '
$ws.Range("D11").Value = '0.07815'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '96.21'
$ws.Range("E12").Value = '  -8.05%  '
$ws.Range("D13").Value = '1.853.81'
$ws.Range("E13").Value = '  -4.32%  '
$ws.Range("D14").Value = '5.109'
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").Value = '0.6621'
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").Value = '279.49'
$ws.Range("E16").Value = '  -5.14%  '
$ws.Range("D17").Value = '30.122.61'
$ws.Range("E17").Value = '  -3.27%  '
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '5.471'
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("E20").Value = '  -3.07%  '
$ws.Range("D21").Value = '2.098.97'
$ws.Range("E21").Value = '  -4.02%  '
$ws.Range("D22").Value = '0.000007211'
$ws.Range("E22").Value = '  -5.29%  '
$ws.Range("D23").Value = '0.9997'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").Value = '6.102'
$ws.Range("E24").Value = '  -5.12%  '
$ws.Range("D25").Value = '9.266'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("D26").Value = '166.28'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("E27").Value = '  -5.18%  '
$ws.Range("D28").Value = '1.895'
$ws.Range("E28").Value = '  -10.70%  '
$ws.Range("D29").Value = '1.330'
$ws.Range("E29").Value = '  -4.65%  '
$ws.Range("D30").Value = '0.09527'
$ws.Range("E30").Value = '  -6.16%  '
$ws.Range("D31").Value = '4.412'
$ws.Range("E31").Value = '  -4.85%  '
$ws.Range("E32").Value = '  -4.50%  '
$ws.Range("D33").Value = '4.075'
$ws.Range("E33").Value = '  -6.54%  '
$ws.Range("D34").Value = '0.04635'
$ws.Range("E34").Value = '  -4.71%  '
$ws.Range("D35").Value = '0.6981'
$ws.Range("E35").Value = '  -6.06%  '
$ws.Range("D36").Value = '1.089'
$ws.Range("E36").Value = '  -3.80%  '
$ws.Range("D37").Value = '2.696'
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("D38").Value = '0.01838'
$ws.Range("E38").Value = '  -6.34%  '
$ws.Range("D39").Value = '6.256'
$ws.Range("E39").Value = '  -4.29%  '
$ws.Range("D40").Value = '2.511'
$ws.Range("E40").Value = '  -4.76%  '
$ws.Range("D41").Value = '72.41'
$ws.Range("E41").Value = '  -6.01%  '
$ws.Range("D42").Value = '0.8523'
$ws.Range("E42").Value = '  -2.50%  '
$ws.Range("E43").Value = '  -6.58%  '
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").Value = '103.69'
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("D46").Value = '0.4128'
$ws.Range("E46").Value = '  -5.51%  '
$ws.Range("D47").Value = '995.97'
$ws.Range("E47").Value = '  -2.76%  '
$ws.Range("D48").Value = '7.157'
$ws.Range("E48").Value = '  -5.53%  '
$ws.Range("D49").Value = '9.166'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").Value = '34.02'
$ws.Range("E50").Value = '  -3.50%  '
$ws.Range("E51").Value = '  -6.59%  '
